$d = $word.ActiveDocument

# 1) Add a new "RAND_SEED" row to the Method-of-Morris options table
#    (the 3rd table in the document: Variable / Type / Values / Description)
$t = $d.Tables.Item(3)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "RAND_SEED"
$newRow.Cells.Item(2).Range.Text = "Unsigned integer"
$newRow.Cells.Item(4).Range.Text = "Seed for the random number generator"

# 2) Drop the stray lastRenderedPageBreak that sits in front of the
#    "The information contained in this file..." sentence. Re-running the
#    paragraph's text through Find/Replace forces Word to rebuild the run,
#    which clears the stale rendering-cache marker.
$findText = "The information contained in this file is described in header line and is pretty much self explanatory "
$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $findText, 2)
